# Insert a new row for table "H004" just above the existing "total" row.
# Before: row 17 = total / 1 / 1  (table_name="total" is the last row)
# After:  row 17 = H004 / 1.1 / 1, and the old "total" row shifts down to row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push the current row 17 ("total", ...) down to make room for the new row.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the H004 data.
$ws.Range("A17").Value = "H004"
$ws.Range("B17").Value = 1.1
$ws.Range("C17").Value = 1
